$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.243.74"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "3.045.74"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.32"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.63"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.440"
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.23"
$ws.Range("E9").Value = "  -4.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.109"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("E11").Value = "  +3.94%  "
$ws.Range("D12").Value = "3.573.17"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.84"
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("D16").Value = "57.101.14"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "3.048.72"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.41"
$ws.Range("E19").Value = "  +4.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.09"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.31"
$ws.Range("E21").Value = "  +2.34%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.507"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.40"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").Value = "3.176.30"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.163"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").Value = "0.0₃0898"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.73"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.16"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.75"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.71"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "152.42"
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("E36").Value = "  +1.38%  "
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.04"
$ws.Range("E38").Value = "  +1.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0670"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").Value = "3.079.19"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.10"
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.662"
$ws.Range("E44").Value = "  +1.96%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.40"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.202.61"
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.959"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.02"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.29"
$ws.Range("E49").Value = "  +4.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0241"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("E51").Value = "  +8.26%  "
